$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.192.96'
$ws.Range('E2').Value = '  +0.15%  '

$ws.Range('D3').Value = '1.601.60'
$ws.Range('E3').Value = '  -0.78%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = '211.94'
$ws.Range('E5').Value = '  -0.61%  '

$ws.Range('E6').Value = '  -0.07%  '

$ws.Range('D7').Value = '0.485'
$ws.Range('E7').Value = '  +0.57%  '

$ws.Range('E8').Value = '  -0.72%  '

$ws.Range('E9').Value = '  -0.83%  '

$ws.Range('D10').Value = '18.12'
$ws.Range('E10').Value = '  -1.70%  '

$ws.Range('D11').Value = '0.0809'
$ws.Range('E11').Value = '  +1.49%  '

$ws.Range('D12').Value = '1.822.91'
$ws.Range('E12').Value = '  -0.85%  '

$ws.Range('D13').Value = '1.607.41'
$ws.Range('E13').Value = '  -0.38%  '

$ws.Range('E14').Value = '  -1.03%  '

$ws.Range('E15').Value = '  +0.87%  '

$ws.Range('D16').Value = '26.182.94'
$ws.Range('E16').Value = '  +0.07%  '

$ws.Range('D17').Value = '61.14'
$ws.Range('E17').Value = '  +0.53%  '

$ws.Range('E18').Value = '  -0.56%  '

$ws.Range('E19').Value = '  -0.01%  '

$ws.Range('D20').Value = '203.62'
$ws.Range('E20').Value = '  +2.53%  '

$ws.Range('E21').Value = '  -0.09%  '

$ws.Range('E22').Value = '  -2.55%  '

$ws.Range('D23').Value = '6.02'
$ws.Range('E23').Value = '  -0.05%  '

$ws.Range('D24').Value = '1.92'
$ws.Range('E24').Value = '  +11.29%  '

$ws.Range('E25').Value = '  +1.43%  '

$ws.Range('E26').Value = '  -0.03%  '

$ws.Range('E27').Value = '  -7.07%  '

$ws.Range('D28').Value = '15.19'
$ws.Range('E28').Value = '  -0.18%  '

$ws.Range('E29').Value = '  +0.03%  '

$ws.Range('D30').Value = '0.0492'
$ws.Range('E30').Value = '  +3.28%  '

$ws.Range('D31').Value = '1.17'
$ws.Range('E31').Value = '  -0.77%  '

$ws.Range('E32').Value = '  -0.33%  '

$ws.Range('D33').Value = '2.91'
$ws.Range('E33').Value = '  -4.36%  '

$ws.Range('E34').Value = '  -2.21%  '

$ws.Range('E35').Value = '  -0.01%  '

$ws.Range('D36').Value = '1.137.54'
$ws.Range('E36').Value = '  +2.65%  '

$ws.Range('E37').Value = '  +6.28%  '

$ws.Range('E38').Value = '  +0.03%  '

$ws.Range('E39').Value = '  -0.55%  '

$ws.Range('D40').Value = '0.784'
$ws.Range('E40').Value = '  -0.94%  '

$ws.Range('E41').Value = '  -2.38%  '

$ws.Range('E42').Value = '  -1.85%  '

$ws.Range('E43').Value = '  +0.82%  '

$ws.Range('D44').Value = '1.737.82'
$ws.Range('E44').Value = '  -0.72%  '

$ws.Range('D45').Value = '92.09'
$ws.Range('E45').Value = '  -1.21%  '

$ws.Range('E46').Value = '  -2.92%  '

$ws.Range('D47').Value = '54.10'
$ws.Range('E47').Value = '  +0.15%  '

$ws.Range('E48').Value = '  -0.75%  '

$ws.Range('D49').Value = '0.406'
$ws.Range('E49').Value = '  -0.78%  '

$ws.Range('E50').Value = '  +0.09%  '

$ws.Range('D51').Value = '0.0₇0947'
$ws.Range('E51').Value = '  -12.81%  '
